# Rename the worksheet from "loginTest" to "loginAsBankManager"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "loginAsBankManager"

# Remove the now-redundant third data row (A3:B3), shrinking the used
# range back down to A1:B2 and dropping the shared-string usage count.
$ws.Rows.Item(3).Delete()
